$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto list data (rows 2-51): Coin (B), Link (C), Price (D), Volume 1h (E).
# Price values that look like plain numbers are prefixed with a literal leading
# apostrophe (via .Formula) so Excel stores them as text, matching the source data
# which uses locale-formatted numbers (e.g. "30.484.76") as plain text strings.
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.484.76', '  +0.59%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.108.42', '  +0.94%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.006', '  +0.58%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '334.48', '  +1.69%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.004', '  +0.54%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5220', '  -0.08%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.4513', '  +4.02%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '53.90', '  +15.69%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08914', '  +0.72%  '),
    @(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.182', '  +1.62%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.14', '  -1.61%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.106.71', '  +0.99%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.816', '  +1.14%  '),
    @(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '8.016', '  +3.45%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '96.69', '  +0.69%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001141', '  +0.82%  '),
    @(18, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.006', '  +0.68%  '),
    @(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06648', '  +0.18%  '),
    @(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '19.21', '  +1.52%  '),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.005', '  +0.60%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.328', '  -0.03%  '),
    @(23, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.549.18', '  +0.61%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '12.43', '  +0.38%  '),
    @(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.350', '  +1.83%  '),
    @(26, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.355.55', '  +1.07%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '22.17', '  -1.11%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.532', '  -2.91%  '),
    @(29, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '162.58', '  +0.43%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '133.82', '  +1.34%  '),
    @(31, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.208', '  +0.01%  '),
    @(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1072', '  +0.09%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.417', '  +3.80%  '),
    @(34, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.635', '  -2.21%  '),
    @(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.944', '  +1.76%  '),
    @(36, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '10.40', '  +4.16%  '),
    @(37, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.759', '  +5.27%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02581', '  -0.20%  '),
    @(39, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06835', '  +2.10%  '),
    @(40, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.2298', '  +1.42%  '),
    @(41, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '12.73', '  +0.35%  '),
    @(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.6868', '  +0.43%  '),
    @(43, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.246', '  -0.09%  '),
    @(44, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.004', '  +0.57%  '),
    @(45, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.318', '  +4.90%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '14.00', '  -0.92%  '),
    @(47, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.6361', '  -0.43%  '),
    @(48, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.662', '  +1.35%  '),
    @(49, 'EOS', 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos', '1.249', '  -0.17%  '),
    @(50, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000349', '  +20.90%  '),
    @(51, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.206', '  +1.32%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $coin = $row[1]
    $link = $row[2]
    $price = $row[3]
    $volume = $row[4]

    $ws.Cells.Item($r, 2).Formula = $coin
    $ws.Cells.Item($r, 3).Formula = $link

    # Force plain-number-looking prices to be stored as text (quote-prefixed),
    # exactly as they already were as inline strings in the source workbook.
    $priceIsNumber = $price -match '^[+-]?[0-9]*\.?[0-9]+$'
    if ($priceIsNumber) {
        $ws.Cells.Item($r, 4).Formula = "'" + $price
    } else {
        $ws.Cells.Item($r, 4).Formula = $price
    }

    $ws.Cells.Item($r, 5).Formula = $volume
}
